# Auto-generated edit script: updates cryptocurrency Price (D) and Volume(1h) (E)
# columns for the Sun Feb 19 22:52:41 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.568.91"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.688.42"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'313.99"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.3886"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").Value = "'0.4029"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "'1.491"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'1.004"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "'53.02"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "'0.08747"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "'25.41"
$ws.Range("E13").Value = "  +7.32%  "
$ws.Range("D14").Value = "'7.513"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").Value = "'0.00001352"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "'7.939"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "1.691.65"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'98.52"
$ws.Range("D19").Value = "'0.07095"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "'7.255"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'14.20"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "24.565.02"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "'2.977"
$ws.Range("E25").Value = "  -8.22%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'22.72"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'161.65"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'8.851"
$ws.Range("E29").Value = "  +16.68%  "
$ws.Range("D30").Value = "'136.88"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'5.229"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "1.879.01"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'0.08810"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").Value = "'7.381"
$ws.Range("E34").Value = "  +4.44%  "
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").Value = "'1.965"
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("D37").Value = "'0.2743"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'0.02912"
$ws.Range("E38").Value = "  +6.80%  "
$ws.Range("E39").Value = "  -5.50%  "
$ws.Range("D40").Value = "'14.23"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("D41").Value = "'0.09121"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "'0.7874"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "'1.455"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "'16.57"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("D45").Value = "'0.7194"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "'2.586"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "'4.197"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'1.340"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "'137.89"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").Value = "'90.93"
$ws.Range("E51").Value = "  -0.09%  "
